$wb = $excel.ActiveWorkbook

# Delete the last two worksheets (pet3_lab1, pet4_lab1)
[void]$wb.Worksheets.Item("pet4_lab1").Delete()
[void]$wb.Worksheets.Item("pet3_lab1").Delete()

# Remove the "keys" row (row 2) from the remaining sheets
for ($i = 1; $i -le 5; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("A2").ClearContents()
}

# Add a new column header on the 5th sheet
$wb.Worksheets.Item(5).Range("B1").Value = "New Column 1"

# Rename the remaining sheets
$wb.Worksheets.Item(1).Name = "pet1_l"
$wb.Worksheets.Item(2).Name = "pet1_a"
$wb.Worksheets.Item(3).Name = "pet1_b"
$wb.Worksheets.Item(4).Name = "pet1_1"
$wb.Worksheets.Item(5).Name = "bip_bop"
